$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to be treated as text so values like "1.009" are not
# auto-converted to numbers by Excel, matching the original inline-string cells.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "26.501.39"
$ws.Range("E2").Value = "  -2.64%  "
$ws.Range("D3").Value = "1.806.31"
$ws.Range("E3").Value = "  -2.50%  "
$ws.Range("E4").Value = "  +0.73%  "
$ws.Range("D5").Value = "1.009"
$ws.Range("E5").Value = "  +0.70%  "
$ws.Range("D6").Value = "308.81"
$ws.Range("E6").Value = "  -1.53%  "
$ws.Range("E7").Value = "  -1.43%  "
$ws.Range("D8").Value = "0.3657"
$ws.Range("E8").Value = "  -1.31%  "
$ws.Range("D9").Value = "0.07118"
$ws.Range("E9").Value = "  -2.29%  "
$ws.Range("D10").Value = "0.8760"
$ws.Range("E10").Value = "  -0.93%  "
$ws.Range("D11").Value = "0.07739"
$ws.Range("E11").Value = "  -1.14%  "
$ws.Range("E12").Value = "  -3.40%  "
$ws.Range("D13").Value = "1.835.33"
$ws.Range("E13").Value = "  -0.46%  "
$ws.Range("D14").Value = "5.266"
$ws.Range("E14").Value = "  -2.06%  "
$ws.Range("E15").Value = "  -2.61%  "
$ws.Range("D16").Value = "86.08"
$ws.Range("D17").Value = "1.012"
$ws.Range("E17").Value = "  +0.88%  "
$ws.Range("D18").Value = "0.000008572"
$ws.Range("E18").Value = "  -3.93%  "
$ws.Range("E19").Value = "  +0.50%  "
$ws.Range("D20").Value = "26.542.56"
$ws.Range("E20").Value = "  -2.56%  "
$ws.Range("D21").Value = "14.24"
$ws.Range("E21").Value = "  -3.20%  "
$ws.Range("E22").Value = "  -2.53%  "
$ws.Range("D23").Value = "10.40"
$ws.Range("E23").Value = "  -1.04%  "
$ws.Range("D24").Value = "1.979"
$ws.Range("E24").Value = "  +2.52%  "
$ws.Range("D25").Value = "150.96"
$ws.Range("E25").Value = "  -0.39%  "
$ws.Range("E26").Value = "  -2.43%  "
$ws.Range("D27").Value = "1.993"
$ws.Range("E27").Value = "  -3.35%  "
$ws.Range("D28").Value = "112.54"
$ws.Range("E28").Value = "  -2.83%  "
$ws.Range("D29").Value = "4.841"
$ws.Range("E29").Value = "  -4.07%  "
$ws.Range("E30").Value = "  -1.78%  "
$ws.Range("D31").Value = "3.036"
$ws.Range("E31").Value = "  -1.86%  "
$ws.Range("D32").Value = "0.7268"
$ws.Range("E32").Value = "  -4.74%  "
$ws.Range("D33").Value = "4.427"
$ws.Range("E33").Value = "  -1.67%  "
$ws.Range("D34").Value = "1.111"
$ws.Range("E34").Value = "  -5.01%  "
$ws.Range("D35").Value = "1.009"
$ws.Range("E35").Value = "  +0.88%  "
$ws.Range("D36").Value = "2.539"
$ws.Range("E36").Value = "  -7.10%  "
$ws.Range("D37").Value = "1.078"
$ws.Range("E37").Value = "  -0.47%  "
$ws.Range("E38").Value = "  -0.95%  "
$ws.Range("E39").Value = "  -2.82%  "
$ws.Range("D40").Value = "2.877"
$ws.Range("E40").Value = "  -2.43%  "
$ws.Range("D41").Value = "6.936"
$ws.Range("E41").Value = "  -1.68%  "
$ws.Range("D42").Value = "0.4990"
$ws.Range("E42").Value = "  -2.15%  "
$ws.Range("D43").Value = "0.1567"
$ws.Range("E43").Value = "  -3.54%  "
$ws.Range("D44").Value = "8.121"
$ws.Range("E44").Value = "  -2.94%  "
$ws.Range("D45").Value = "1.010"
$ws.Range("E45").Value = "  +0.83%  "
$ws.Range("D46").Value = "0.4608"
$ws.Range("E46").Value = "  -3.76%  "
$ws.Range("D47").Value = "101.82"
$ws.Range("E47").Value = "  -0.26%  "
$ws.Range("D48").Value = "9.930"
$ws.Range("E48").Value = "  -3.69%  "
$ws.Range("E49").Value = "  -3.12%  "
$ws.Range("D50").Value = "0.06001"
$ws.Range("E50").Value = "  -3.34%  "
$ws.Range("D51").Value = "63.84"
$ws.Range("E51").Value = "  -2.49%  "

# Restore default (Normal) style on column D so no stray number-format styling
# is left behind on cells, matching the original unstyled cells.
$ws.Range("D2:D51").Style = "Normal"

